$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B149").Value = 48654
$ws.Range("E149").Value = 38.26
$ws.Range("F149").Value = -1
$ws.Range("G149").Value = -32.02
$ws.Range("B150").Value = 63902
$ws.Range("E150").Value = 34.04
$ws.Range("F150").Value = 2
$ws.Range("G150").Value = 64.04000000000001
$ws.Range("B183").Value = 57552
$ws.Range("E183").Value = 136.86
$ws.Range("F183").Value = -5
$ws.Range("G183").Value = -603.45
$ws.Range("B184").Value = 64329
$ws.Range("E184").Value = 128.32
$ws.Range("F184").Value = 6
$ws.Range("G184").Value = 724.14
$ws.Range("B313").Value = 62997
$ws.Range("F313").Value = 72
$ws.Range("G313").Value = 22020.48
$ws.Range("B314").Value = 57854
$ws.Range("F314").Value = 2
$ws.Range("G314").Value = 611.6799999999999
$ws.Range("B316").Value = 57077
$ws.Range("D316").Value = 93.08
$ws.Range("E316").Value = 111.2
$ws.Range("F316").Value = 1
$ws.Range("G316").Value = 93.08
$ws.Range("B317").Value = 61610
$ws.Range("E317").Value = 122.71
$ws.Range("F317").Value = -58
$ws.Range("G317").Value = -5957.18
$ws.Range("B318").Value = 63565
$ws.Range("D318").Value = 102.71
$ws.Range("E318").Value = 109.19
$ws.Range("F318").Value = 60
$ws.Range("G318").Value = 6162.6
$ws.Range("B346").Value = 63520
$ws.Range("E346").Value = 153.4
$ws.Range("F346").Value = 97
$ws.Range("G346").Value = 13995.16
$ws.Range("B347").Value = 55373
$ws.Range("E347").Value = 163.62
$ws.Range("F347").Value = -94
$ws.Range("G347").Value = -13562.32
$ws.Range("B351").Value = 63531
$ws.Range("E351").Value = 152.53
$ws.Range("F351").Value = 80
$ws.Range("G351").Value = 11478.4
$ws.Range("B352").Value = 57802
$ws.Range("E352").Value = 162.71
$ws.Range("F352").Value = -79
$ws.Range("G352").Value = -11334.92
$ws.Range("B372").Value = 63652
$ws.Range("E372").Value = 55.42
$ws.Range("F372").Value = 250
$ws.Range("G372").Value = 13032.5
$ws.Range("B373").Value = 57885
$ws.Range("E373").Value = 62.28
$ws.Range("F373").Value = 4
$ws.Range("G373").Value = 208.52
$ws.Range("B375").Value = 63563
$ws.Range("E375").Value = 119.04
$ws.Range("F375").Value = 15
$ws.Range("G375").Value = 1679.4
$ws.Range("B376").Value = 61605
$ws.Range("E376").Value = 133.78
$ws.Range("F376").Value = -13
$ws.Range("G376").Value = -1455.48
$ws.Range("B379").Value = 61608
$ws.Range("E379").Value = 154.12
$ws.Range("F379").Value = -56
$ws.Range("G379").Value = -7224.56
$ws.Range("B380").Value = 63564
$ws.Range("E380").Value = 137.16
$ws.Range("F380").Value = 57
$ws.Range("G380").Value = 7353.57
$ws.Range("B382").Value = 60325
$ws.Range("E382").Value = 151.57
$ws.Range("F382").Value = -102
$ws.Range("G382").Value = -12939.72
$ws.Range("B383").Value = 63560
$ws.Range("E383").Value = 134.87
$ws.Range("F383").Value = 104
$ws.Range("G383").Value = 13193.44
$ws.Range("B389").Value = 62865
$ws.Range("F389").Value = 151
$ws.Range("G389").Value = 12051.31
$ws.Range("B390").Value = 57817
$ws.Range("F390").Value = 3
$ws.Range("G390").Value = 239.43
$ws.Range("B400").Value = 62933
$ws.Range("F400").Value = 146
$ws.Range("G400").Value = 8632.98
$ws.Range("B401").Value = 57835
$ws.Range("F401").Value = 1
$ws.Range("G401").Value = 59.13
$ws.Range("B419").Value = 63007
$ws.Range("F419").Value = 984
$ws.Range("G419").Value = 168588.72
$ws.Range("B420").Value = 57856
$ws.Range("F420").Value = 2
$ws.Range("G420").Value = 342.66
$ws.Range("B421").Value = 63008
$ws.Range("F421").Value = 504
$ws.Range("G421").Value = 76189.67999999999
$ws.Range("B422").Value = 57857
$ws.Range("F422").Value = 3
$ws.Range("G422").Value = 453.51
$ws.Range("B431").Value = 63102
$ws.Range("C431").Value = "HUL-Vim Bar Multipack Fw 4X200G"
$ws.Range("F431").Value = 36
$ws.Range("G431").Value = 2140.92
$ws.Range("B432").Value = 53082
$ws.Range("C432").Value = "HUL-VIM BAR MULTIPACK FW 4X200G"
$ws.Range("F432").Value = 1
$ws.Range("G432").Value = 59.47
$ws.Range("B536").Value = 47097
$ws.Range("D536").Value = 112.28
$ws.Range("E536").Value = 134.16
$ws.Range("F536").Value = 15
$ws.Range("G536").Value = 1684.2
$ws.Range("B537").Value = 58047
$ws.Range("D537").Value = 105.54
$ws.Range("E537").Value = 126.1
$ws.Range("F537").Value = 54
$ws.Range("G537").Value = 5699.16
$ws.Range("B579").Value = 65069
$ws.Range("E579").Value = 14.3
$ws.Range("F579").Value = 172
$ws.Range("G579").Value = 2313.4
$ws.Range("B580").Value = 53757
$ws.Range("E580").Value = 16.08
$ws.Range("F580").Value = -159
$ws.Range("G580").Value = -2138.55
$ws.Range("B583").Value = 65066
$ws.Range("E583").Value = 13.61
$ws.Range("F583").Value = 313
$ws.Range("G583").Value = 4009.53
$ws.Range("B584").Value = 53263
$ws.Range("E584").Value = 15.29
$ws.Range("F584").Value = -309
$ws.Range("G584").Value = -3958.29
$ws.Range("B590").Value = 45706
$ws.Range("E590").Value = 23.58
$ws.Range("F590").Value = -202
$ws.Range("G590").Value = -3985.46
$ws.Range("B591").Value = 64922
$ws.Range("E591").Value = 20.98
$ws.Range("F591").Value = 207
$ws.Range("G591").Value = 4084.11
$ws.Range("B593").Value = 64927
$ws.Range("E593").Value = 17.26
$ws.Range("F593").Value = 295
$ws.Range("G593").Value = 4784.9
$ws.Range("B594").Value = 45718
$ws.Range("E594").Value = 19.38
$ws.Range("F594").Value = -294
$ws.Range("G594").Value = -4768.68
$ws.Range("B599").Value = 64925
$ws.Range("E599").Value = 13.97
$ws.Range("F599").Value = 302
$ws.Range("G599").Value = 3971.3
$ws.Range("B600").Value = 45709
$ws.Range("E600").Value = 15.69
$ws.Range("F600").Value = -300
$ws.Range("G600").Value = -3945
$ws.Range("B601").Value = 45702
$ws.Range("E601").Value = 31.43
$ws.Range("F601").Value = -215
$ws.Range("G601").Value = -5654.5
$ws.Range("B602").Value = 64919
$ws.Range("E602").Value = 27.97
$ws.Range("F602").Value = 224
$ws.Range("G602").Value = 5891.2
$ws.Range("B687").Value = 53319
$ws.Range("E687").Value = 310.64
$ws.Range("F687").Value = -6
$ws.Range("G687").Value = -1643.52
$ws.Range("B688").Value = 64810
$ws.Range("E688").Value = 291.22
$ws.Range("F688").Value = 7
$ws.Range("G688").Value = 1917.44
$ws.Range("B709").Value = 60025
$ws.Range("E709").Value = 37.22
$ws.Range("F709").Value = -98
$ws.Range("G709").Value = -3217.34
$ws.Range("B710").Value = 64833
$ws.Range("E710").Value = 34.9
$ws.Range("F710").Value = 99
$ws.Range("G710").Value = 3250.17
$ws.Range("B720").Value = 60022
$ws.Range("E720").Value = 37.22
$ws.Range("F720").Value = -113
$ws.Range("G720").Value = -3709.79
$ws.Range("B721").Value = 64830
$ws.Range("E721").Value = 34.9
$ws.Range("F721").Value = 117
$ws.Range("G721").Value = 3841.11
$ws.Range("B872").Value = 65079
$ws.Range("F872").Value = 21
$ws.Range("G872").Value = 858.27
$ws.Range("B873").Value = 65362
$ws.Range("F873").Value = 2
$ws.Range("G873").Value = 81.73999999999999
